# Actualización automática 2025-06-16 13:01:14
#
# The "VENTAS POR GRUPO" sheet gains a new category column ("GRANITO")
# inserted before the "GRIFERIAS" column, plus three new category columns
# appended at the end ("NO RESURTIBLES", "PANELES PVC", "PANELES PU").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

$lastDataRow = 54
$totalsRow = 55

# --- Insert the GRANITO column before column F (GRIFERIAS), shifting the
#     rest of the category columns one slot to the right. Excel carries the
#     formatting of the surrounding cells into the freshly inserted column,
#     same as interactively right-clicking "Insert" on the column header. ---
$ws.Columns("F:F").Insert()

$ws.Cells.Item(1, 6).Value = "GRANITO"
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# --- Append three brand-new category columns after the last existing one
#     (old column N "SAL SOLUBLE", now column O after the insert above).
#     Clone formatting from column O into P:R, then fill in the data. ---
$ws.Range("O1:O" + $totalsRow).Copy()
$ws.Range("P1:R" + $totalsRow).PasteSpecial(-4122)

$ws.Cells.Item(1, 16).Value = "NO RESURTIBLES"
$ws.Cells.Item(1, 17).Value = "PANELES PVC"
$ws.Cells.Item(1, 18).Value = "PANELES PU"

for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
    $ws.Cells.Item($r, 18).Value = 0
}

$ws.Cells.Item($totalsRow, 16).Value = "0 de 53"
$ws.Cells.Item($totalsRow, 17).Value = "0 de 53"
$ws.Cells.Item($totalsRow, 18).Value = "0 de 53"

# --- Column widths: GRANITO matches the narrow numeric columns, the three
#     appended columns get their own widths. (ColumnWidth is in character
#     units; Excel's internal storage adds ~5/6 of a character vs. the
#     OOXML <col width> value for the default Calibri 11 font, so we back
#     that out here to land on the exact target widths.) ---
$charOffset = 5 / 6
$ws.Columns.Item(6).ColumnWidth = 13 - $charOffset
$ws.Columns.Item(16).ColumnWidth = 20 - $charOffset
$ws.Columns.Item(17).ColumnWidth = 17 - $charOffset
$ws.Columns.Item(18).ColumnWidth = 16 - $charOffset

Write-Output "done"
